$wb = $excel.ActiveWorkbook

# Rename several sheet tabs
$wb.Worksheets.Item("total mortality").Name = "mortality rates"
$wb.Worksheets.Item("mortality").Name = "causes of death"
$wb.Worksheets.Item("RRStunting").Name = "RR death by stunting"
$wb.Worksheets.Item("RRWasting").Name = "RR death by wasting"
$wb.Worksheets.Item("RRBreastfeeding").Name = "RR death by breastfeeding"
$wb.Worksheets.Item("RR Death by Birth Outcome").Name = "RR death by birth outcome"
$wb.Worksheets.Item("OR stunting for complements").Name = "OR stunting by compfeeding"
$wb.Worksheets.Item("OR appropriateBF by interv").Name = "OR correctBF by interventn"

# Update the column title "Complements group" -> "Food security & education"
$ws = $wb.Worksheets.Item("OR stunting by compfeeding")
$ws.Range("A1").Value = "Food security & education"

# Make "mortality rates" (formerly "total mortality") the active sheet/tab
$wb.Worksheets.Item("mortality rates").Activate()
